$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new rows 9-14 with data (Name, Instansi, NIP_NUPTK, Tema Webinar)
$ws.Range("A9").Value = "Adit"
$ws.Range("B9").Value = "TP"
$ws.Range("C9").Value = 901
$ws.Range("D9").Value = "Test"

$ws.Range("A10").Value = "Iqbal"
$ws.Range("B10").Value = "IPS"
$ws.Range("C10").Value = 902
$ws.Range("D10").Value = "Test"

$ws.Range("A11").Value = "Aldi"
$ws.Range("B11").Value = "IPA"
$ws.Range("C11").Value = 903
$ws.Range("D11").Value = "Test"

$ws.Range("A12").Value = "Bagas"
$ws.Range("B12").Value = "DKV"
$ws.Range("C12").Value = 904
$ws.Range("D12").Value = "Test"

$ws.Range("A13").Value = "Basium"
$ws.Range("B13").Value = "YTBR"
$ws.Range("C13").Value = 905
$ws.Range("D13").Value = "Test"

$ws.Range("A14").Value = "Tegar"
$ws.Range("B14").Value = "BTK"
$ws.Range("C14").Value = 906
$ws.Range("D14").Value = "Test"

# Match the existing centered style used by the rest of the table (style index 1)
$ws.Range("A9:D14").HorizontalAlignment = -4108

# Rows 15-17 pick up the centered style in column D (still empty)
$ws.Range("D15:D17").HorizontalAlignment = -4108

# Set selection to G7, matching the final cursor position
$ws.Range("G7").Select()
